$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text formatting so numeric-looking strings (e.g. "52.469.88", "2.00") are preserved as text

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '52.469.88'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.90%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.022.01'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.10%  '

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '356.41'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.81%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '110.31'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.35%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.566'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.49%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.06%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.627'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -1.10%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.75'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.45%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0864'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -3.52%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '19.46'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.29%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.482.42'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.64%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.84'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -2.76%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.008.23'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.60%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +3.10%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '52.526.03'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.56'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +8.43%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.59'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.16%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.78'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.65%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.0₃0979'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -1.32%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '69.83'

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.09'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.27%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.76'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.85%  '

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.15%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '27.25'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.27%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.70'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +2.85%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.13%  '

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.57%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.47'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +1.78%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '10.41'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -3.04%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '37.13'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.93%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.20'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +17.51%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '50.99'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -4.05%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.35%  '

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.26'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -4.57%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.22'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -3.87%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.00'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -2.95%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.73'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.48%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.44%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '23.24'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -2.19%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '123.63'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +7.90%  '

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.29%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.137.28'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -1.66%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.43'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -4.23%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.40'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.250'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.39%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0339'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.84%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.943'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.50%  '
